$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns.
# D-column values are plain text (e.g. "55.715.30", "0.999") and must stay
# text even when they look numeric, so force a text format before writing,
# then restore the default "Normal" style so formatting matches the source.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '55.715.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.972.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '491.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.47%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.970.59'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.420'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E11").Value = '  -4.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.352'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.47%  '
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.478.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '55.538.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.957.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000141'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '322.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.31%  '
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.466'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '60.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -11.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0852'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0653'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.000.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.26'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.632'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.135.52'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0237'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.61'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.25%  '
